$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.674.48"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.852.06"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4239"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07302"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "1.853.81"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06895"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008926"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").Value = "27.688.03"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.994"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.62%  "
$ws.Range("D25").Value = "2.078.38"
$ws.Range("E25").Value = "  -5.64%  "
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.275"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.885"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08867"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7688"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.565"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.976"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05360"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01937"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.816"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.904"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.313"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06548"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.629"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.49%  "
